$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Part 1 - Table caption paragraph:
#   "Table {SEQ Table \* ARABIC}: List and definition of symbols..."
#   becomes
#   "Table 1: List and definition of [_GoBack]symbols..."
# (the SEQ field is flattened to plain text, and the editing-cursor
#  "_GoBack" bookmark moves here from its old spot in the table body)
# ---------------------------------------------------------------------

# Remove the SEQ field; this collapses "Table " + field + ": List..." down
# to "Table : List and definition of symbols used in the text and
# equations with their units." while preserving the _Ref bookmark that
# wraps "Table ".
$f = $d.Fields.Item(1)
$f.Delete()

$full = $d.Content.Text
$tableWordStart = $full.IndexOf("Table")
$colonPos = $full.IndexOf(":", $tableWordStart)

# "Table " (with trailing space) -> "Table" (no trailing space)
$rTable = $d.Range($tableWordStart, $colonPos)
$rTable.Text = "Table"
$afterTableEnd = $tableWordStart + 5

# Insert " 1" right after "Table"
$rInsertPos = $d.Range($afterTableEnd, $afterTableEnd)
$rInsertPos.InsertAfter(" 1")

# Force "Table" / " 1" to remain distinct runs: toggle formatting on
# "Table" only, then release it via a freshly-fetched Range (avoids the
# engine re-coalescing the two pieces back into a single run).
$rTableOnly = $d.Range($tableWordStart, $afterTableEnd)
$rTableOnly.Font.Bold = 1
$rTableOnlyAgain = $d.Range($tableWordStart, $afterTableEnd)
$rTableOnlyAgain.Font.Bold = 0

# Split ": List and definition of symbols used in the text and equations
# with their units." into ": List and definition of " + "symbols used...".
$full = $d.Content.Text
$afterOneEnd = $afterTableEnd + 2
$symbolsPos = $full.IndexOf("symbols used")
$rLead = $d.Range($afterOneEnd, $symbolsPos)
$rLead.Font.Bold = 1
$rLeadAgain = $d.Range($afterOneEnd, $symbolsPos)
$rLeadAgain.Font.Bold = 0

# Drop a fresh "_GoBack" bookmark at the split point.
$rGoBack = $d.Range($symbolsPos, $symbolsPos)
$d.Bookmarks.Add("_GoBack", $rGoBack)

# ---------------------------------------------------------------------
# Part 2 - table body paragraph:
#   "... scaling of anabo[_GoBack]lism and energy expenditure"
#   becomes
#   "... scaling of anabolism and energy expenditure"
# (drops the stray _GoBack bookmark that used to sit mid-word and
#  reunites the run that it had split in two)
# ---------------------------------------------------------------------
$d.Content.Find.Execute("anabolism", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "anabolism", 2)

Write-Output "Caption now reads: $($d.Paragraphs.Item(1).Range.Text)"
